$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit A: "404 page" -> "404-page" and drop the two proofErr markers
# (gramStart/gramEnd) that wrap it, while keeping the surrounding runs
# ("Handle invalid user id/", " (4) ", "X") intact as separate runs.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Handle invalid user id/404 page ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $startPos = $rng.Start
    $endPos = $rng.End

    $part1 = "Handle invalid user id/"
    $part2 = "404-page"
    $part3 = " "
    $rest = "(4) "
    # " (4) " and "X" already follow in the document and are left as-is
    # text-wise; only the matched span ("Handle invalid user id/404 page ",
    # which ends at the leading space of the " (4) " run) is replaced.
    # $rest is only used here to compute where "X" starts -- it is NOT
    # (re)written, since it is already present right after the edit range.

    # Replacing the whole span (run1 + proofErr + run2 + proofErr + leading
    # space of run3) collapses it to one run and drops the now-stale
    # proofErr markers that bracketed "404 page".
    $editRange = $d.Range($startPos, $endPos)
    $editRange.Text = $part1 + $part2 + $part3

    # Re-introduce the original run boundaries by dropping a bookmark at
    # each split point and immediately deleting it again -- Bookmarks.Add
    # splits the run it lands in without touching (merging) any of the
    # other runs in the paragraph.
    $split1 = $startPos + $part1.Length
    $split2 = $split1 + $part2.Length
    $split3 = $split2 + $part3.Length + $rest.Length

    $d.Bookmarks.Add("ZZEditSplit1", $d.Range($split1, $split1))
    $d.Bookmarks("ZZEditSplit1").Delete()
    $d.Bookmarks.Add("ZZEditSplit2", $d.Range($split2, $split2))
    $d.Bookmarks("ZZEditSplit2").Delete()
    $d.Bookmarks.Add("ZZEditSplit3", $d.Range($split3, $split3))
    $d.Bookmarks("ZZEditSplit3").Delete()
}

# ---------------------------------------------------------------------
# Edit B: move the "_GoBack" bookmark from just before "Bug: " down into
# the middle of "without" ("withou" | "t sidebar") in the next bullet.
# Bookmarks.Add with an existing name re-seats the single bookmark, and
# (like above) splits only the run it lands in.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("with and withou", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $rng2.End
    $d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))
}
